$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.440.78"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.291.66"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'301.08"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").Value = "'95.06"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").Value = "'0.506"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").Value = "'34.35"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("D11").Value = "'19.00"
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "'6.71"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "2.645.35"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "2.296.28"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "'0.779"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "42.369.31"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "'12.17"
$ws.Range("E19").Value = "  -6.53%  "
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "'5.96"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "'2.26"
$ws.Range("E23").Value = "  +6.34%  "
$ws.Range("D24").Value = "'235.53"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("D27").Value = "'24.22"
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("D28").Value = "'2.21"
$ws.Range("E28").Value = "  -7.23%  "
$ws.Range("D29").Value = "'164.34"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("D30").Value = "'9.03"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").Value = "'31.64"
$ws.Range("E31").Value = "  -4.21%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "'4.97"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "'17.47"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").Value = "'0.0695"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  -2.96%  "
$ws.Range("D37").Value = "'4.35"
$ws.Range("E37").Value = "  -8.56%  "
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").Value = "'2.68"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").Value = "'19.87"
$ws.Range("E42").Value = "  +8.15%  "
$ws.Range("D43").Value = "1.947.55"
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("D44").Value = "'10.33"
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").Value = "'2.09"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.514.81"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "'52.83"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "'2.80"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'70.51"
$ws.Range("E51").Value = "  -1.38%  "
